$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Arabic "name" value for the MLE/ara row (B5): it incorrectly held
# the Arabic word for "female" (أنثى) and should hold the word for "male"
# (الذكر). Setting the value adds a new shared string and repoints B5 at it.
$ws.Range("B5").Value = "الذكر"

# Autofit column B (name) to its (new) longest entry, as Excel would do
# after the content changed, producing a <cols> entry for column B.
$ws.Columns("B").AutoFit()

# Configure the page setup (paper size + orientation) so a <pageSetup>
# element is emitted for the sheet.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Leave the active selection on D16, matching the saved view state.
$ws.Range("D16").Select()
